# 1st changes of mifos to finflux
#
# The "Repayment schedule" sheet gains a new blank column between the
# existing "Late" column (was N) and its neighbours (old N/O/P shift to
# O/P/Q). The sheet also becomes the active tab/sheet, with the
# selection left on cell R6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (moves tabSelected/activeTab here).
$ws.Activate()

# Insert a new blank column at N - shifts old N (Late), O (Outstanding)
# and P (Disbursement) one column to the right (-> O, P, Q).
$ws.Columns("N").Insert()

# The inherited column width for the new column matches the width of
# column M (10.71 characters) that sits immediately to its left.
$ws.Columns("N").ColumnWidth = 9.83

# Leave the final selection on R6, which also marks this sheet/tab as
# the active one.
$ws.Range("R6").Select()
